$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that needs to move
# forward by 2 days (46070 -> 46072) for each data row (rows 2 through 8).
for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 2
}
